# Updated symbol list on Sat Jan 14 16:32:10 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coinranking.com crypto snapshot on sheet1. Values are written as literal
# text (matching the existing inline-string cells) rather than numbers, so
# a leading apostrophe forces text entry and ClearFormats() strips the
# resulting quote-prefix styling back to the sheet's default look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}


Set-TextValue "D2" "303.79"
Set-TextValue "E2" "5.78%"
Set-TextValue "D3" "32.46"
Set-TextValue "E3" "11.17%"
Set-TextValue "D4" "5.293"
Set-TextValue "E4" "1.97%"
Set-TextValue "D5" "0.07499"
Set-TextValue "E5" "7.80%"
Set-TextValue "D6" "7.804"
Set-TextValue "E6" "5.54%"
Set-TextValue "D7" "3.800"
Set-TextValue "E7" "6.88%"
Set-TextValue "D8" "1.513"
Set-TextValue "E8" "7.22%"
Set-TextValue "D9" "0.9189"
Set-TextValue "E9" "2.17%"
Set-TextValue "D10" "0.01770"
Set-TextValue "E10" "2,630.70%"
Set-TextValue "D11" "0.1697"
Set-TextValue "E11" "6.25%"
Set-TextValue "D12" "0.07797"
Set-TextValue "E12" "6.07%"
Set-TextValue "D13" "0.08066"
Set-TextValue "E13" "4.72%"
Set-TextValue "D14" "0.03043"
Set-TextValue "E14" "4.08%"
Set-TextValue "D15" "0.09904"
Set-TextValue "E15" "10.22%"
Set-TextValue "D16" "0.001492"
Set-TextValue "E16" "-6.50%"
Set-TextValue "D17" "0.04601"
Set-TextValue "E17" "1.58%"
Set-TextValue "D18" "0.006527"
Set-TextValue "E18" "-1.38%"
Set-TextValue "D19" "3.479"
Set-TextValue "E19" "0.46%"
Set-TextValue "D20" "2.228"
Set-TextValue "E20" "0.15%"
Set-TextValue "E21" "3.86%"
Set-TextValue "D22" "0.1333"
Set-TextValue "E22" "0.67%"
Set-TextValue "D23" "4.565"
Set-TextValue "E23" "13.43%"
Set-TextValue "D24" "0.1619"
Set-TextValue "E24" "3.86%"
Set-TextValue "D25" "0.001218"
Set-TextValue "E25" "0.59%"
Set-TextValue "D26" "0.004441"
Set-TextValue "E26" "1.81%"
Set-TextValue "D27" "0.0001400"
Set-TextValue "E27" "19.69%"
Set-TextValue "D28" "0.0001738"
Set-TextValue "E28" "7.42%"
Set-TextValue "E40" "5.78%"
Set-TextValue "D41" "0.007214"
Set-TextValue "E41" "3.71%"
Set-TextValue "D42" "0.1346"
Set-TextValue "E42" "8.40%"
Set-TextValue "D43" "0.002171"
Set-TextValue "E43" "4.36%"
Set-TextValue "D44" "0.01268"
Set-TextValue "E44" "7.99%"
Set-TextValue "D45" "0.00006031"
Set-TextValue "E45" "3.54%"
Set-TextValue "D46" "0.7082"
Set-TextValue "E46" "-63.29%"
Set-TextValue "D47" "0.01298"
Set-TextValue "E47" "-0.71%"
